$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.413.63'
$ws.Range('E2').Value = '  +0.96%  '
$ws.Range('D3').Value = '2.327.20'
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''302.45'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').Value = '''98.12'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = '''0.504'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('D10').Value = '''35.69'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('D11').Value = '''19.55'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +7.61%  '
$ws.Range('D12').Value = '''0.0799'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').Value = '''6.92'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('D15').Value = '2.696.15'
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('D16').Value = '2.342.47'
$ws.Range('E16').Value = '  +1.66%  '
$ws.Range('D17').Value = '''0.793'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '43.322.91'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').Value = '''12.82'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('D20').Value = '0.0₃0901'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').Value = '''6.08'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('D22').Value = '''68.07'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = '''237.41'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').Value = '''2.25'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.71%  '
$ws.Range('D25').Value = '''2.46'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '''25.05'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('D28').Value = '''2.21'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +7.31%  '
$ws.Range('D29').Value = '''164.49'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('D30').Value = '''9.14'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('D31').Value = '''33.27'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').Value = '''17.85'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.72%  '
$ws.Range('D35').Value = '''4.48'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -7.30%  '
$ws.Range('D36').Value = '''0.0705'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.05%  '
$ws.Range('E37').Value = '  -1.17%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('D42').Value = '1.988.75'
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '''19.33'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +8.51%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''10.64'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +5.65%  '
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '''2.07'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').Value = '''2.81'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.561.00'
$ws.Range('E48').Value = '  +1.12%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '''54.08'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.60%  '
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('D51').Value = '''72.75'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.75%  '
